$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "68.674.70"
$ws.Range("E2").Value = "  +0.54%  "

Set-TextValue "D3" "2.712.73"
$ws.Range("E3").Value = "  +2.47%  "

Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue "D5" "599.87"
$ws.Range("E5").Value = "  +0.38%  "

Set-TextValue "D6" "160.84"
$ws.Range("E6").Value = "  +3.00%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.23%  "

Set-TextValue "D9" "2.712.99"
$ws.Range("E9").Value = "  +2.54%  "

Set-TextValue "D10" "0.140"
$ws.Range("E10").Value = "  -3.98%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("E13").Value = "  +1.89%  "

$ws.Range("E14").Value = "  +1.10%  "

Set-TextValue "D15" "3.209.37"
$ws.Range("E15").Value = "  +2.58%  "

$ws.Range("E16").Value = "  -2.54%  "

Set-TextValue "D17" "68.707.64"
$ws.Range("E17").Value = "  +0.70%  "

Set-TextValue "D18" "2.712.68"
$ws.Range("E18").Value = "  +2.89%  "

Set-TextValue "D19" "11.85"
$ws.Range("E19").Value = "  +4.13%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "367.07"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "7.65"
$ws.Range("E21").Value = "  +2.42%  "

$ws.Range("E22").Value = "  +3.07%  "

Set-TextValue "D23" "4.94"
$ws.Range("E23").Value = "  +2.04%  "

Set-TextValue "D24" "2.12"
$ws.Range("E24").Value = "  +2.90%  "

Set-TextValue "D25" "75.32"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -0.02%  "

Set-TextValue "D27" "10.16"
$ws.Range("E27").Value = "  +4.89%  "

Set-TextValue "D28" "2.818.77"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("E29").Value = "  -1.07%  "

Set-TextValue "D30" "584.47"
$ws.Range("E30").Value = "  +4.40%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  +3.29%  "

$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("E34").Value = "  +5.07%  "

Set-TextValue "D35" "1.65"
$ws.Range("E35").Value = "  +5.52%  "

$ws.Range("E36").Value = "  +1.72%  "

Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  +0.08%  "

Set-TextValue "D38" "20.28"
$ws.Range("E38").Value = "  +4.82%  "

Set-TextValue "D39" "161.37"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  +2.47%  "

Set-TextValue "D41" "1.90"
$ws.Range("E41").Value = "  +0.78%  "

Set-TextValue "D42" "5.44"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.66"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D44" "17.86"
$ws.Range("E44").Value = "  +0.47%  "

Set-TextValue "D45" "0.0₆0320"
$ws.Range("E45").Value = "  -5.70%  "

$ws.Range("E46").Value = "  +0.07%  "

Set-TextValue "D47" "160.10"
$ws.Range("E47").Value = "  +0.25%  "

Set-TextValue "D48" "3.93"
$ws.Range("E48").Value = "  +5.01%  "

Set-TextValue "D49" "1.78"
$ws.Range("E49").Value = "  +5.46%  "

Set-TextValue "D50" "0.608"
$ws.Range("E50").Value = "  +8.34%  "

Set-TextValue "D51" "22.42"
$ws.Range("E51").Value = "  +2.00%  "
